$d = $word.ActiveDocument

$d.Content.Find.Execute("787÷3=262, 1", $true, $false, $false, $false, $false, $true, 1, $false, "906÷9=100, 6", 2) | Out-Null
$d.Content.Find.Execute("356÷6=59, 2", $true, $false, $false, $false, $false, $true, 1, $false, "389÷9=43, 2", 2) | Out-Null
$d.Content.Find.Execute("342÷2=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "931÷3=310, 1", 2) | Out-Null
$d.Content.Find.Execute("934÷6=155, 4", $true, $false, $false, $false, $false, $true, 1, $false, "576÷2=288, 0", 2) | Out-Null
$d.Content.Find.Execute("364÷2=182, 0", $true, $false, $false, $false, $false, $true, 1, $false, "516÷7=73, 5", 2) | Out-Null
$d.Content.Find.Execute("392÷4=98, 0", $true, $false, $false, $false, $false, $true, 1, $false, "215÷3=71, 2", 2) | Out-Null
$d.Content.Find.Execute("783÷6=130, 3", $true, $false, $false, $false, $false, $true, 1, $false, "405÷2=202, 1", 2) | Out-Null
$d.Content.Find.Execute("295÷2=147, 1", $true, $false, $false, $false, $false, $true, 1, $false, "306÷3=102, 0", 2) | Out-Null
$d.Content.Find.Execute("691÷8=86, 3", $true, $false, $false, $false, $false, $true, 1, $false, "998÷7=142, 4", 2) | Out-Null
$d.Content.Find.Execute("370÷3=123, 1", $true, $false, $false, $false, $false, $true, 1, $false, "596÷4=149, 0", 2) | Out-Null
$d.Content.Find.Execute("910÷4=227, 2", $true, $false, $false, $false, $false, $true, 1, $false, "438÷8=54, 6", 2) | Out-Null
$d.Content.Find.Execute("328÷5=65, 3", $true, $false, $false, $false, $false, $true, 1, $false, "163÷3=54, 1", 2) | Out-Null
$d.Content.Find.Execute("444÷6=74, 0", $true, $false, $false, $false, $false, $true, 1, $false, "968÷2=484, 0", 2) | Out-Null
$d.Content.Find.Execute("872÷3=290, 2", $true, $false, $false, $false, $false, $true, 1, $false, "314÷7=44, 6", 2) | Out-Null
$d.Content.Find.Execute("873÷8=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "103÷4=25, 3", 2) | Out-Null
$d.Content.Find.Execute("139÷2=69, 1", $true, $false, $false, $false, $false, $true, 1, $false, "442÷7=63, 1", 2) | Out-Null
$d.Content.Find.Execute("569÷7=81, 2", $true, $false, $false, $false, $false, $true, 1, $false, "982÷7=140, 2", 2) | Out-Null
$d.Content.Find.Execute("758÷8=94, 6", $true, $false, $false, $false, $false, $true, 1, $false, "943÷8=117, 7", 2) | Out-Null
$d.Content.Find.Execute("183÷3=61, 0", $true, $false, $false, $false, $false, $true, 1, $false, "494÷6=82, 2", 2) | Out-Null
$d.Content.Find.Execute("795÷8=99, 3", $true, $false, $false, $false, $false, $true, 1, $false, "883÷4=220, 3", 2) | Out-Null
$d.Content.Find.Execute("789÷8=98, 5", $true, $false, $false, $false, $false, $true, 1, $false, "254÷3=84, 2", 2) | Out-Null
$d.Content.Find.Execute("806÷9=89, 5", $true, $false, $false, $false, $false, $true, 1, $false, "615÷8=76, 7", 2) | Out-Null
$d.Content.Find.Execute("478÷2=239, 0", $true, $false, $false, $false, $false, $true, 1, $false, "757÷9=84, 1", 2) | Out-Null
$d.Content.Find.Execute("727÷2=363, 1", $true, $false, $false, $false, $false, $true, 1, $false, "618÷2=309, 0", 2) | Out-Null
$d.Content.Find.Execute("474÷2=237, 0", $true, $false, $false, $false, $false, $true, 1, $false, "238÷3=79, 1", 2) | Out-Null
